$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 195, pushing the existing data (old rows 195-210)
# down to rows 196-211. This also extends the sheet dimension and carries
# the date-column formatting down onto the new row.
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new weekly record.
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44461
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 100112032
$ws.Range("G195").Value = "Zapallo italiano"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 90
$ws.Range("K195").Value = 11000
$ws.Range("L195").Value = 11500
$ws.Range("M195").Value = 11278
$ws.Range("N195").Value = "`$/caja 70 unidades"
$ws.Range("O195").Value = "Región de Arica y Parinacota"
$ws.Range("P195").Value = 161
$ws.Range("Q195").Value = 70
$ws.Range("R195").Value = "Hortaliza"
